# "login associated files modify 1"
# Sheet1 contains a small "User" table in columns G:H with field names in H1:H3.
# The field that used to be labelled "Password" is being renamed to "user_pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Password" label (H3) to "user_pass"
$ws.Range("H3").Value = "user_pass"

# Widen column A (matches the author's manual column resize)
$ws.Columns.Item(1).ColumnWidth = 12.8

# Move the active selection to E7 (last edited cell)
$ws.Range("E7").Select()
